$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "covid19_cases_switzerland": fill in row 19 (new daily data)
# and update the AB19 running-total formula.
# -----------------------------------------------------------------
$wsCases = $wb.Worksheets.Item("covid19_cases_switzerland")

$wsCases.Range("B19").Value = 234
$wsCases.Range("C19").Value = 4
$wsCases.Range("D19").Value = 27
$wsCases.Range("F19").Value = 292
$wsCases.Range("G19").Value = 376
$wsCases.Range("H19").Value = 189
$wsCases.Range("J19").Value = 29
$wsCases.Range("P19").Value = 21
$wsCases.Range("S19").Value = 83
$wsCases.Range("T19").Value = 73
$wsCases.Range("X19").Value = 1880
$wsCases.Range("Y19").Value = 441

$wsCases.Range("AB19").Formula = "=SUM(B19:AA19)+Z16+O18+N18+L18+I18"

# -----------------------------------------------------------------
# Sheet "Quellen": selection moved to B13
# -----------------------------------------------------------------
$wsQuellen = $wb.Worksheets.Item("Quellen")
$wsQuellen.Activate()
$wsQuellen.Range("B13").Select()

# -----------------------------------------------------------------
# Sheet "demographics": updated incidence figures (column M) from BAG
# -----------------------------------------------------------------
$wsDemo = $wb.Worksheets.Item("demographics")

$wsDemo.Range("M2").Value = 62.8
$wsDemo.Range("M3").Value = 43
$wsDemo.Range("M4").Value = 235.3
$wsDemo.Range("M5").Value = 34.5
$wsDemo.Range("M6").Value = 36.8
$wsDemo.Range("M7").Value = 136.3
$wsDemo.Range("M8").Value = 34.7
$wsDemo.Range("M9").Value = 326.9
$wsDemo.Range("M10").Value = 128.2
$wsDemo.Range("M11").Value = 59.3
$wsDemo.Range("M12").Value = 101.3
$wsDemo.Range("M13").Value = 23.9
$wsDemo.Range("M14").Value = 30.4
$wsDemo.Range("M15").Value = 133.6
$wsDemo.Range("M16").Value = 222.8
$wsDemo.Range("M17").Value = 84.3
$wsDemo.Range("M18").Value = 45.9
$wsDemo.Range("M19").Value = 32.3
$wsDemo.Range("M20").Value = 36.6
$wsDemo.Range("M21").Value = 68.1
$wsDemo.Range("M22").Value = 48.9
$wsDemo.Range("M23").Value = 78.7
$wsDemo.Range("M24").Value = 71.8
$wsDemo.Range("M25").Value = 55.5
$wsDemo.Range("M26").Value = 49.4
$wsDemo.Range("M27").Value = 24.8

# -----------------------------------------------------------------
# Re-activate the "covid19_cases_switzerland" sheet and move its
# selection last, so it ends up as the active tab (as in the source).
# -----------------------------------------------------------------
$wsCases.Activate()
$wsCases.Range("AB24").Select()
